$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Materias primas" (raw materials) column values to reorder the
# ingredient lists, per the shared-strings diff.
$ws.Range("C2").Value = "leche,harina,huevos,vainilla,"
$ws.Range("C3").Value = "manzana,huevos,harina,"
$ws.Range("C4").Value = "vainilla,harina,huevos,"
$ws.Range("C6").Value = "merengue,harina,crema, limon, huevos,"
